$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$tr.Text = "Per il Server si genera una tabella, come da specifica, con valori parzialmente casuali.`rOutput *classifica_giudici_1_svc(void *in, struct svc_req *rqstp): per ogni giudice viene calcolato il punteggio in base ai partecipanti e successivamente viene ordinato il risultato, restituendo al Client l’array ordinato di Giudici.`rint *esprimi_voto_1_svc(Input *input, struct svc_req *rqstp): si verifica l’esistenza del partecipante passato dal Client e si aggiunge o sottrae un voto, restituendo al Client un valore diverso da -1.`r`r`r"

$tr.Characters(1, 7).Font.Size = 1600

$tr.Characters(8, 6).Font.Size = 1600
$tr.Characters(8, 6).Font.Color.RGB = 0xF0B000

$tr.Characters(14, 75).Font.Size = 1600

$tr.Characters(90, 6).Font.Size = 1600
$tr.Characters(90, 6).Font.Name = "Consolas"
$tr.Characters(90, 6).Font.Color.ObjectThemeColor = 10

$tr.Characters(96, 2).Font.Size = 1600
$tr.Characters(96, 2).Font.Name = "Consolas"
$tr.Characters(96, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(98, 24).Font.Size = 1600
$tr.Characters(98, 24).Font.Name = "Consolas"
$tr.Characters(98, 24).Font.Color.RGB = 0xAADCDC

$tr.Characters(122, 1).Font.Size = 1600
$tr.Characters(122, 1).Font.Name = "Consolas"
$tr.Characters(122, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(123, 4).Font.Size = 1600
$tr.Characters(123, 4).Font.Name = "Consolas"
$tr.Characters(123, 4).Font.Color.RGB = 0xD69C56

$tr.Characters(127, 2).Font.Size = 1600
$tr.Characters(127, 2).Font.Name = "Consolas"
$tr.Characters(127, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(129, 2).Font.Size = 1600
$tr.Characters(129, 2).Font.Name = "Consolas"
$tr.Characters(129, 2).Font.Color.RGB = 0xFEDC9C

$tr.Characters(131, 2).Font.Size = 1600
$tr.Characters(131, 2).Font.Name = "Consolas"
$tr.Characters(131, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(133, 6).Font.Size = 1600
$tr.Characters(133, 6).Font.Name = "Consolas"
$tr.Characters(133, 6).Font.Color.RGB = 0xD69C56

$tr.Characters(139, 1).Font.Size = 1600
$tr.Characters(139, 1).Font.Name = "Consolas"
$tr.Characters(139, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(140, 7).Font.Size = 1600
$tr.Characters(140, 7).Font.Name = "Consolas"
$tr.Characters(140, 7).Font.Color.ObjectThemeColor = 10

$tr.Characters(147, 2).Font.Size = 1600
$tr.Characters(147, 2).Font.Name = "Consolas"
$tr.Characters(147, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(149, 5).Font.Size = 1600
$tr.Characters(149, 5).Font.Name = "Consolas"
$tr.Characters(149, 5).Font.Color.RGB = 0xFEDC9C

$tr.Characters(154, 1).Font.Size = 1600
$tr.Characters(154, 1).Font.Name = "Consolas"
$tr.Characters(154, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(155, 169).Font.Size = 1600

$tr.Characters(325, 3).Font.Size = 1600
$tr.Characters(325, 3).Font.Name = "Consolas"
$tr.Characters(325, 3).Font.Color.RGB = 0xD69C56

$tr.Characters(328, 2).Font.Size = 1600
$tr.Characters(328, 2).Font.Name = "Consolas"
$tr.Characters(328, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(330, 18).Font.Size = 1600
$tr.Characters(330, 18).Font.Name = "Consolas"
$tr.Characters(330, 18).Font.Color.RGB = 0xAADCDC

$tr.Characters(348, 1).Font.Size = 1600
$tr.Characters(348, 1).Font.Name = "Consolas"
$tr.Characters(348, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(349, 5).Font.Size = 1600
$tr.Characters(349, 5).Font.Name = "Consolas"
$tr.Characters(349, 5).Font.Color.ObjectThemeColor = 10

$tr.Characters(354, 2).Font.Size = 1600
$tr.Characters(354, 2).Font.Name = "Consolas"
$tr.Characters(354, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(356, 5).Font.Size = 1600
$tr.Characters(356, 5).Font.Name = "Consolas"
$tr.Characters(356, 5).Font.Color.RGB = 0xFEDC9C

$tr.Characters(361, 2).Font.Size = 1600
$tr.Characters(361, 2).Font.Name = "Consolas"
$tr.Characters(361, 2).Font.Color.RGB = 0xD4D4D4

$tr.Characters(363, 6).Font.Size = 1600
$tr.Characters(363, 6).Font.Name = "Consolas"
$tr.Characters(363, 6).Font.Color.RGB = 0xD69C56

$tr.Characters(369, 1).Font.Size = 1600
$tr.Characters(369, 1).Font.Name = "Consolas"
$tr.Characters(369, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(370, 7).Font.Size = 1600
$tr.Characters(370, 7).Font.Name = "Consolas"
$tr.Characters(370, 7).Font.Color.ObjectThemeColor = 10

$tr.Characters(377, 1).Font.Size = 1600
$tr.Characters(377, 1).Font.Name = "Consolas"
$tr.Characters(377, 1).Font.Color.ObjectThemeColor = 10

$tr.Characters(378, 1).Font.Size = 1600
$tr.Characters(378, 1).Font.Name = "Consolas"
$tr.Characters(378, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(379, 5).Font.Size = 1600
$tr.Characters(379, 5).Font.Name = "Consolas"
$tr.Characters(379, 5).Font.Color.RGB = 0xFEDC9C

$tr.Characters(384, 1).Font.Size = 1600
$tr.Characters(384, 1).Font.Name = "Consolas"
$tr.Characters(384, 1).Font.Color.RGB = 0xD4D4D4

$tr.Characters(385, 67).Font.Size = 1600

$tr.Characters(452, 19).Font.Size = 1600

$tr.Characters(471, 55).Font.Size = 1600

$sh.Height = 243.51007874015747

